$wb = $excel.ActiveWorkbook

# --- Sheet handles -----------------------------------------------------
$wsPopDef   = $wb.Worksheets.Item("Population Definitions")
$wsXferDef  = $wb.Worksheets.Item("Transfer Definitions")
$wsXferDet  = $wb.Worksheets.Item("Transfer Details")
$wsPopSize  = $wb.Worksheets.Item("Population Sizes")
$wsPrev     = $wb.Worksheets.Item("Prevalence")
$wsMort     = $wb.Worksheets.Item("Mortality Rates")
$wsEpi      = $wb.Worksheets.Item("Epidemic Characteristics")
$wsCascade  = $wb.Worksheets.Item("Cascade Parameters")

# --- Transfer Details: bump the default migration fraction to a Number -
# D2 was "Fraction" -> "Number", and the companion formula's default
# changes from 0.1 to 10 to match.
$wsXferDet.Range("D2").Value = "Number"
$wsXferDet.Range("E2").Formula = '=IF(A2<>"...",IF(SUMPRODUCT(--(G2:V2<>""))=0,10,"N.A."),"")'

# --- Prevalence: add a "Recovered" compartment block (rows 25-27) ------
$wsPrev.Range("A25").Value = "Recovered"
$wsPrev.Range("B25").Value = "Format"
$wsPrev.Range("C25").Value = "Assumption"
$years = 2000..2015
for ($i = 0; $i -lt $years.Count; $i++) {
    $col = 5 + $i   # column E = 5
    $wsPrev.Cells.Item(25, $col).Value = $years[$i]
}

$wsPrev.Range("A26").Formula = "='Population Definitions'!`$A`$2"
$wsPrev.Range("B26").Value = "Number"
$wsPrev.Range("C26").Value = 2500
$wsPrev.Range("D26").Value = "OR"

$wsPrev.Range("A27").Formula = "='Population Definitions'!`$A`$3"
$wsPrev.Range("B27").Value = "Number"
$wsPrev.Range("C27").Value = 2500
$wsPrev.Range("D27").Value = "OR"

# Extend the "Fraction,Number" dropdown validation onto the new rows too
$wsPrev.Range("B26").Validation.Add(3, 1, 1, "Fraction,Number")
$wsPrev.Range("B27").Validation.Add(3, 1, 1, "Fraction,Number")

# --- View-state bookkeeping ---------------------------------------------
# Cascade Parameters: scroll back to the top (drops topLeftCell="A6")
[void]$wsCascade.Activate()
$excel.ActiveWindow.ScrollRow = 1
[void]$wsCascade.Range("C43").Select()

# Transfer Definitions: selection moves to B31 (tab no longer active)
[void]$wsXferDef.Activate()
[void]$wsXferDef.Range("B31").Select()

# Transfer Details: selection moves to E3
[void]$wsXferDet.Activate()
[void]$wsXferDet.Range("E3").Select()

# Population Sizes: selection moves to C5
[void]$wsPopSize.Activate()
[void]$wsPopSize.Range("C5").Select()

# Prevalence: becomes the active tab, scrolled to row 7, selection C28
[void]$wsPrev.Activate()
$excel.ActiveWindow.ScrollRow = 7
[void]$wsPrev.Range("C28").Select()

Write-Output "done"
